$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.320.05'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '3.132.92'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.128.63'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.150'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.472'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000255'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.32'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '3.655.77'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('E16').Value = '  +3.28%  '
$ws.Range('D17').Value = '64.238.82'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '3.139.23'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '474.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.56'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.718'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.96%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.117'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.15%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.61'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.76%  '
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.89'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.47'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.04%  '
$ws.Range('D38').Value = '0.0₃0733'
$ws.Range('E38').Value = '  +4.32%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '449.82'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.96'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.80%  '
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.118'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('D44').Value = '2.844.80'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.265'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.54%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '26.15'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.113'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.31%  '
